# Quiz table: choose columns for math questions/answers and
# question/answer pairs. Builds a 4-column header + data table on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: labels identifying each column's role for the quiz picker.
$ws.Range("A1").Value = "m1#MATH"
$ws.Range("B1").Value = "q1#QUESTION"
$ws.Range("C1").Value = "q2#QUESTION"
$ws.Range("D1").Value = "m2#MATH"

# Row 2
$ws.Range("A2").Value = "1 + 2 = 3"
$ws.Range("C2").Value = "can you read? -> Yes."
$ws.Range("D2").Value = "2+5 = 7"

# Row 3
$ws.Range("A3").Value = "2 / 5 = .4"
$ws.Range("C3").Value = "who uses the computer? -> me."
$ws.Range("D3").Value = "1 ** 0 = undef"

# Row 4
$ws.Range("D4").Value = "lim(x->5, x-5) = 0"

# A2 picked up an explicit (default-valued) protection lock, as if the
# author opened Format Cells > Protection on that cell.
$ws.Range("A2").Locked = $true
